# Update "NV-30 Đào Vương Anh 8-2024.xlsx":
#  1) Insert a new worksheet "Đơn phụ phẫu 1" between "Đơn sale chính" and "Lương"
#     containing the phụ phẫu 1 orders detail table.
#  2) Refresh the computed totals on the "Lương" sheet that are driven by that
#     new data (công / tiền phụ phẫu 1 tại LONG XUYÊN, and the resulting totals).

$wb = $excel.ActiveWorkbook

$saleSheet = $wb.Worksheets.Item(1)

# --- 1) Insert the new sheet right after "Đơn sale chính" -------------------
$newSheet = $wb.Worksheets.Add($null, $saleSheet)
$newSheet.Name = "Đơn phụ phẫu 1"

# Header row
$newSheet.Range("A1").Value = "Tiền tố"
$newSheet.Range("B1").Value = "Mã dịch vụ"
$newSheet.Range("C1").Value = "Ngày thực hiện"
$newSheet.Range("D1").Value = "Cơ sở"
$newSheet.Range("E1").Value = "Khách hàng"
$newSheet.Range("F1").Value = "Nguồn khách"
$newSheet.Range("G1").Value = "Tên dịch vụ"
$newSheet.Range("H1").Value = "Phụ phẫu 1"
$newSheet.Range("I1").Value = "Công phụ phẫu 1"

# Row 2
$newSheet.Range("A2").Value = "HD-LUXURY"
$newSheet.Range("B2").Value = 616
$newSheet.Range("C2").Value = "'08-02-2024"
$newSheet.Range("D2").Value = "LONG XUYÊN"
$newSheet.Range("E2").Value = "Chị duyên"
$newSheet.Range("F2").Value = "Khách cũ giới thiệu"
$newSheet.Range("G2").Value = "Cắt mí"
$newSheet.Range("H2").Value = "Đào Vương Anh"
$newSheet.Range("I2").Value = 50000

# Row 3
$newSheet.Range("A3").Value = "HD-LUXURY"
$newSheet.Range("B3").Value = 617
$newSheet.Range("C3").Value = "'08-02-2024"
$newSheet.Range("D3").Value = "LONG XUYÊN"
$newSheet.Range("E3").Value = "Cô tú"
$newSheet.Range("F3").Value = "Khách cũ"
$newSheet.Range("G3").Value = "Nâng cung chân mày"
$newSheet.Range("H3").Value = "Đào Vương Anh"
$newSheet.Range("I3").Value = 50000

# Row 4 - totals
$newSheet.Range("A4").Value = "Tổng"
$newSheet.Range("B4").Value = 2
$newSheet.Range("I4").Value = 100000

# --- 2) Refresh "Lương" sheet totals affected by the new data ---------------
# Re-fetch the sheet reference *after* the insertion, since the new sheet
# shifted sheet positions.
$luongSheet = $wb.Worksheets.Item("Lương")
$luongSheet.Range("B12").Value = 2
$luongSheet.Range("B13").Value = 70000
$luongSheet.Range("B14").Value = 285714.2857142857
$luongSheet.Range("B19").Value = 100000
$luongSheet.Range("B33").Value = 455714.2857142857
$luongSheet.Range("A35").Value = "Tổng lương tại HỆ THỐNG"
$luongSheet.Range("B35").Value = 455714.2857142857
